$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 7 now has real dates: update the header label ---
$ws.Range("H1").Value = "Week 7 -- June 13 - 19"

# --- Fill in this week's (Week 7) scores for athletes who competed ---
$ws.Range("H2").Value = 306.1
$ws.Range("H3").Value = 203.1
$ws.Range("H4").Value = 216.1
$ws.Range("H5").Value = 208.2
$ws.Range("H6").Value = 269.6
$ws.Range("H7").Value = 480
$ws.Range("H8").Value = 189.6
$ws.Range("H11").Value = 429

# --- Widen column H now that it holds a full "Week 7 -- June 13 - 19" label ---
$ws.Columns.Item(8).ColumnWidth = 18.3

# --- A few new styled-but-empty rows appear further down column B ---
$ws.Range("B23").NumberFormat = "0.0"
$ws.Range("B30").NumberFormat = "0.0"
$ws.Range("B37").NumberFormat = "0.0"

# --- Move the active selection ---
$ws.Range("D14").Select() | Out-Null
